$wb = $excel.ActiveWorkbook
$g = [double]"5.686312626471138e+23"

$ws = $wb.Worksheets.Item("ROW35-FE-LIFTER")
$dateFmt = $ws.Cells.Item(2,1).NumberFormat
$ws.Cells.Item(65,1).Value = 45721.72894113426
$ws.Cells.Item(65,1).NumberFormat = $dateFmt
$ws.Cells.Item(65,2).Value = "0x01,0x90"
$ws.Cells.Item(65,3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x14,0x41,0x0c,"
$ws.Cells.Item(65,4).Value = "0x01,0x90,"
$ws.Cells.Item(65,5).Value = "0xd"
$ws.Cells.Item(65,6).Value = 400
$ws.Cells.Item(65,7).Value = $g
$ws.Cells.Item(65,8).Value = 400
$ws.Cells.Item(65,9).Value = 13
$ws.Cells.Item(66,1).Value = 45721.72896428241
$ws.Cells.Item(66,1).NumberFormat = $dateFmt
$ws.Cells.Item(66,2).Value = "0x01,0x90"
$ws.Cells.Item(66,3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x14,0x41,0x0c,"
$ws.Cells.Item(66,4).Value = "0x01,0x90,"
$ws.Cells.Item(66,5).Value = "0xd"
$ws.Cells.Item(66,6).Value = 400
$ws.Cells.Item(66,7).Value = $g
$ws.Cells.Item(66,8).Value = 400
$ws.Cells.Item(66,9).Value = 13
$ws.Cells.Item(67,1).Value = 45721.72898766203
$ws.Cells.Item(67,1).NumberFormat = $dateFmt
$ws.Cells.Item(67,2).Value = "0x01,0x90"
$ws.Cells.Item(67,3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x14,0x41,0x0c,"
$ws.Cells.Item(67,4).Value = "0x01,0x90,"
$ws.Cells.Item(67,5).Value = "0xd"
$ws.Cells.Item(67,6).Value = 400
$ws.Cells.Item(67,7).Value = $g
$ws.Cells.Item(67,8).Value = 400
$ws.Cells.Item(67,9).Value = 13
$ws.Cells.Item(68,1).Value = 45722.22908445602
$ws.Cells.Item(68,1).NumberFormat = $dateFmt
$ws.Cells.Item(68,2).Value = "0x01,0x90"
$ws.Cells.Item(68,3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x14,0x41,0x0c,"
$ws.Cells.Item(68,4).Value = "0x01,0x90,"
$ws.Cells.Item(68,5).Value = "0xd"
$ws.Cells.Item(68,6).Value = 400
$ws.Cells.Item(68,7).Value = $g
$ws.Cells.Item(68,8).Value = 400
$ws.Cells.Item(68,9).Value = 13
$ws.Cells.Item(69,1).Value = 45722.22910648148
$ws.Cells.Item(69,1).NumberFormat = $dateFmt
$ws.Cells.Item(69,2).Value = "0x01,0x90"
$ws.Cells.Item(69,3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x14,0x41,0x0c,"
$ws.Cells.Item(69,4).Value = "0x01,0x90,"
$ws.Cells.Item(69,5).Value = "0xd"
$ws.Cells.Item(69,6).Value = 400
$ws.Cells.Item(69,7).Value = $g
$ws.Cells.Item(69,8).Value = 400
$ws.Cells.Item(69,9).Value = 13
$ws.Cells.Item(70,1).Value = 45722.22912973379
$ws.Cells.Item(70,1).NumberFormat = $dateFmt
$ws.Cells.Item(70,2).Value = "0x01,0x90"
$ws.Cells.Item(70,3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x14,0x41,0x0c,"
$ws.Cells.Item(70,4).Value = "0x01,0x90,"
$ws.Cells.Item(70,5).Value = "0xd"
$ws.Cells.Item(70,6).Value = 400
$ws.Cells.Item(70,7).Value = $g
$ws.Cells.Item(70,8).Value = 400
$ws.Cells.Item(70,9).Value = 13
$ws.Cells.Item(71,1).Value = 45723.19113143518
$ws.Cells.Item(71,1).NumberFormat = $dateFmt
$ws.Cells.Item(71,2).Value = "0x01,0x90"
$ws.Cells.Item(71,3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x14,0x41,0x0c,"
$ws.Cells.Item(71,4).Value = "0x01,0x90,"
$ws.Cells.Item(71,5).Value = "0xd"
$ws.Cells.Item(71,6).Value = 400
$ws.Cells.Item(71,7).Value = $g
$ws.Cells.Item(71,8).Value = 400
$ws.Cells.Item(71,9).Value = 13
$ws.Cells.Item(72,1).Value = 45723.19115481481
$ws.Cells.Item(72,1).NumberFormat = $dateFmt
$ws.Cells.Item(72,2).Value = "0x01,0x90"
$ws.Cells.Item(72,3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x14,0x41,0x0c,"
$ws.Cells.Item(72,4).Value = "0x01,0x90,"
$ws.Cells.Item(72,5).Value = "0xd"
$ws.Cells.Item(72,6).Value = 400
$ws.Cells.Item(72,7).Value = $g
$ws.Cells.Item(72,8).Value = 400
$ws.Cells.Item(72,9).Value = 13
$ws.Cells.Item(73,1).Value = 45723.19117797454
$ws.Cells.Item(73,1).NumberFormat = $dateFmt
$ws.Cells.Item(73,2).Value = "0x01,0x90"
$ws.Cells.Item(73,3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x14,0x41,0x0c,"
$ws.Cells.Item(73,4).Value = "0x01,0x90,"
$ws.Cells.Item(73,5).Value = "0xd"
$ws.Cells.Item(73,6).Value = 400
$ws.Cells.Item(73,7).Value = $g
$ws.Cells.Item(73,8).Value = 400
$ws.Cells.Item(73,9).Value = 13

$ws = $wb.Worksheets.Item("ROW35-MID-LIFTER")
$dateFmt = $ws.Cells.Item(2,1).NumberFormat
$ws.Cells.Item(68,1).Value = 45721.72988806713
$ws.Cells.Item(68,1).NumberFormat = $dateFmt
$ws.Cells.Item(68,2).Value = "0x01,0x90"
$ws.Cells.Item(68,3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c,"
$ws.Cells.Item(68,4).Value = "0x01,0x90,"
$ws.Cells.Item(68,5).Value = "0xe"
$ws.Cells.Item(68,6).Value = 400
$ws.Cells.Item(68,7).Value = $g
$ws.Cells.Item(68,8).Value = 400
$ws.Cells.Item(68,9).Value = 14
$ws.Cells.Item(69,1).Value = 45721.72991133102
$ws.Cells.Item(69,1).NumberFormat = $dateFmt
$ws.Cells.Item(69,2).Value = "0x01,0x90"
$ws.Cells.Item(69,3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c,"
$ws.Cells.Item(69,4).Value = "0x01,0x90,"
$ws.Cells.Item(69,5).Value = "0xe"
$ws.Cells.Item(69,6).Value = 400
$ws.Cells.Item(69,7).Value = $g
$ws.Cells.Item(69,8).Value = 400
$ws.Cells.Item(69,9).Value = 14
$ws.Cells.Item(70,1).Value = 45721.72993467592
$ws.Cells.Item(70,1).NumberFormat = $dateFmt
$ws.Cells.Item(70,2).Value = "0x01,0x90"
$ws.Cells.Item(70,3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c,"
$ws.Cells.Item(70,4).Value = "0x01,0x90,"
$ws.Cells.Item(70,5).Value = "0xe"
$ws.Cells.Item(70,6).Value = 400
$ws.Cells.Item(70,7).Value = $g
$ws.Cells.Item(70,8).Value = 400
$ws.Cells.Item(70,9).Value = 14
$ws.Cells.Item(71,1).Value = 45722.23020512731
$ws.Cells.Item(71,1).NumberFormat = $dateFmt
$ws.Cells.Item(71,2).Value = "0x01,0x90"
$ws.Cells.Item(71,3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c,"
$ws.Cells.Item(71,4).Value = "0x01,0x90,"
$ws.Cells.Item(71,5).Value = "0xe"
$ws.Cells.Item(71,6).Value = 400
$ws.Cells.Item(71,7).Value = $g
$ws.Cells.Item(71,8).Value = 400
$ws.Cells.Item(71,9).Value = 14
$ws.Cells.Item(72,1).Value = 45722.23022724537
$ws.Cells.Item(72,1).NumberFormat = $dateFmt
$ws.Cells.Item(72,2).Value = "0x01,0x90"
$ws.Cells.Item(72,3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c,"
$ws.Cells.Item(72,4).Value = "0x01,0x90,"
$ws.Cells.Item(72,5).Value = "0xe"
$ws.Cells.Item(72,6).Value = 400
$ws.Cells.Item(72,7).Value = $g
$ws.Cells.Item(72,8).Value = 400
$ws.Cells.Item(72,9).Value = 14
$ws.Cells.Item(73,1).Value = 45722.23025050926
$ws.Cells.Item(73,1).NumberFormat = $dateFmt
$ws.Cells.Item(73,2).Value = "0x01,0x90"
$ws.Cells.Item(73,3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c,"
$ws.Cells.Item(73,4).Value = "0x01,0x90,"
$ws.Cells.Item(73,5).Value = "0xe"
$ws.Cells.Item(73,6).Value = 400
$ws.Cells.Item(73,7).Value = $g
$ws.Cells.Item(73,8).Value = 400
$ws.Cells.Item(73,9).Value = 14
$ws.Cells.Item(74,1).Value = 45723.19124234954
$ws.Cells.Item(74,1).NumberFormat = $dateFmt
$ws.Cells.Item(74,2).Value = "0x01,0x90"
$ws.Cells.Item(74,3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c,"
$ws.Cells.Item(74,4).Value = "0x01,0x90,"
$ws.Cells.Item(74,5).Value = "0xe"
$ws.Cells.Item(74,6).Value = 400
$ws.Cells.Item(74,7).Value = $g
$ws.Cells.Item(74,8).Value = 400
$ws.Cells.Item(74,9).Value = 14
$ws.Cells.Item(75,1).Value = 45723.19126548611
$ws.Cells.Item(75,1).NumberFormat = $dateFmt
$ws.Cells.Item(75,2).Value = "0x01,0x90"
$ws.Cells.Item(75,3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c,"
$ws.Cells.Item(75,4).Value = "0x01,0x90,"
$ws.Cells.Item(75,5).Value = "0xe"
$ws.Cells.Item(75,6).Value = 400
$ws.Cells.Item(75,7).Value = $g
$ws.Cells.Item(75,8).Value = 400
$ws.Cells.Item(75,9).Value = 14
$ws.Cells.Item(76,1).Value = 45723.19128880787
$ws.Cells.Item(76,1).NumberFormat = $dateFmt
$ws.Cells.Item(76,2).Value = "0x01,0x90"
$ws.Cells.Item(76,3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c,"
$ws.Cells.Item(76,4).Value = "0x01,0x90,"
$ws.Cells.Item(76,5).Value = "0xe"
$ws.Cells.Item(76,6).Value = 400
$ws.Cells.Item(76,7).Value = $g
$ws.Cells.Item(76,8).Value = 400
$ws.Cells.Item(76,9).Value = 14

$ws = $wb.Worksheets.Item("ROW02-FE-LIFTER")
$dateFmt = $ws.Cells.Item(2,1).NumberFormat
$ws.Cells.Item(65,1).Value = 45721.72979079861
$ws.Cells.Item(65,1).NumberFormat = $dateFmt
$ws.Cells.Item(65,2).Value = "0x01,0x90"
$ws.Cells.Item(65,3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x06,0x41,0x0c,"
$ws.Cells.Item(65,4).Value = "0x01,0x90,"
$ws.Cells.Item(65,5).Value = "0xff"
$ws.Cells.Item(65,6).Value = 400
$ws.Cells.Item(65,7).Value = $g
$ws.Cells.Item(65,8).Value = 400
$ws.Cells.Item(65,9).Value = 255
$ws.Cells.Item(66,1).Value = 45721.72981402778
$ws.Cells.Item(66,1).NumberFormat = $dateFmt
$ws.Cells.Item(66,2).Value = "0x01,0x90"
$ws.Cells.Item(66,3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x06,0x41,0x0c,"
$ws.Cells.Item(66,4).Value = "0x01,0x90,"
$ws.Cells.Item(66,5).Value = "0xff"
$ws.Cells.Item(66,6).Value = 400
$ws.Cells.Item(66,7).Value = $g
$ws.Cells.Item(66,8).Value = 400
$ws.Cells.Item(66,9).Value = 255
$ws.Cells.Item(67,1).Value = 45721.7298375463
$ws.Cells.Item(67,1).NumberFormat = $dateFmt
$ws.Cells.Item(67,2).Value = "0x01,0x90"
$ws.Cells.Item(67,3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x06,0x41,0x0c,"
$ws.Cells.Item(67,4).Value = "0x01,0x90,"
$ws.Cells.Item(67,5).Value = "0xff"
$ws.Cells.Item(67,6).Value = 400
$ws.Cells.Item(67,7).Value = $g
$ws.Cells.Item(67,8).Value = 400
$ws.Cells.Item(67,9).Value = 255
$ws.Cells.Item(68,1).Value = 45722.23010768519
$ws.Cells.Item(68,1).NumberFormat = $dateFmt
$ws.Cells.Item(68,2).Value = "0x01,0x90"
$ws.Cells.Item(68,3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x06,0x41,0x0c,"
$ws.Cells.Item(68,4).Value = "0x01,0x90,"
$ws.Cells.Item(68,5).Value = "0xff"
$ws.Cells.Item(68,6).Value = 400
$ws.Cells.Item(68,7).Value = $g
$ws.Cells.Item(68,8).Value = 400
$ws.Cells.Item(68,9).Value = 255
$ws.Cells.Item(69,1).Value = 45722.23012966435
$ws.Cells.Item(69,1).NumberFormat = $dateFmt
$ws.Cells.Item(69,2).Value = "0x01,0x90"
$ws.Cells.Item(69,3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x06,0x41,0x0c,"
$ws.Cells.Item(69,4).Value = "0x01,0x90,"
$ws.Cells.Item(69,5).Value = "0xff"
$ws.Cells.Item(69,6).Value = 400
$ws.Cells.Item(69,7).Value = $g
$ws.Cells.Item(69,8).Value = 400
$ws.Cells.Item(69,9).Value = 255
$ws.Cells.Item(70,1).Value = 45722.23015302084
$ws.Cells.Item(70,1).NumberFormat = $dateFmt
$ws.Cells.Item(70,2).Value = "0x01,0x90"
$ws.Cells.Item(70,3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x06,0x41,0x0c,"
$ws.Cells.Item(70,4).Value = "0x01,0x90,"
$ws.Cells.Item(70,5).Value = "0xff"
$ws.Cells.Item(70,6).Value = 400
$ws.Cells.Item(70,7).Value = $g
$ws.Cells.Item(70,8).Value = 400
$ws.Cells.Item(70,9).Value = 255
$ws.Cells.Item(71,1).Value = 45723.19132489583
$ws.Cells.Item(71,1).NumberFormat = $dateFmt
$ws.Cells.Item(71,2).Value = "0x01,0x90"
$ws.Cells.Item(71,3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x06,0x41,0x0c,"
$ws.Cells.Item(71,4).Value = "0x01,0x90,"
$ws.Cells.Item(71,5).Value = "0xff"
$ws.Cells.Item(71,6).Value = 400
$ws.Cells.Item(71,7).Value = $g
$ws.Cells.Item(71,8).Value = 400
$ws.Cells.Item(71,9).Value = 255
$ws.Cells.Item(72,1).Value = 45723.19134847222
$ws.Cells.Item(72,1).NumberFormat = $dateFmt
$ws.Cells.Item(72,2).Value = "0x01,0x90"
$ws.Cells.Item(72,3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x06,0x41,0x0c,"
$ws.Cells.Item(72,4).Value = "0x01,0x90,"
$ws.Cells.Item(72,5).Value = "0xff"
$ws.Cells.Item(72,6).Value = 400
$ws.Cells.Item(72,7).Value = $g
$ws.Cells.Item(72,8).Value = 400
$ws.Cells.Item(72,9).Value = 255
$ws.Cells.Item(73,1).Value = 45723.19137149306
$ws.Cells.Item(73,1).NumberFormat = $dateFmt
$ws.Cells.Item(73,2).Value = "0x01,0x90"
$ws.Cells.Item(73,3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x06,0x41,0x0c,"
$ws.Cells.Item(73,4).Value = "0x01,0x90,"
$ws.Cells.Item(73,5).Value = "0xff"
$ws.Cells.Item(73,6).Value = 400
$ws.Cells.Item(73,7).Value = $g
$ws.Cells.Item(73,8).Value = 400
$ws.Cells.Item(73,9).Value = 255

$ws = $wb.Worksheets.Item("ROW02-MID-LIFTER")
$dateFmt = $ws.Cells.Item(2,1).NumberFormat
$ws.Cells.Item(65,1).Value = 45721.72820228009
$ws.Cells.Item(65,1).NumberFormat = $dateFmt
$ws.Cells.Item(65,2).Value = "0x01,0x90"
$ws.Cells.Item(65,3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
$ws.Cells.Item(65,4).Value = "0x01,0x90,"
$ws.Cells.Item(65,5).Value = "0x3"
$ws.Cells.Item(65,6).Value = 400
$ws.Cells.Item(65,7).Value = $g
$ws.Cells.Item(65,8).Value = 400
$ws.Cells.Item(65,9).Value = 3
$ws.Cells.Item(66,1).Value = 45721.7282258449
$ws.Cells.Item(66,1).NumberFormat = $dateFmt
$ws.Cells.Item(66,2).Value = "0x01,0x90"
$ws.Cells.Item(66,3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
$ws.Cells.Item(66,4).Value = "0x01,0x90,"
$ws.Cells.Item(66,5).Value = "0x3"
$ws.Cells.Item(66,6).Value = 400
$ws.Cells.Item(66,7).Value = $g
$ws.Cells.Item(66,8).Value = 400
$ws.Cells.Item(66,9).Value = 3
$ws.Cells.Item(67,1).Value = 45721.72824888889
$ws.Cells.Item(67,1).NumberFormat = $dateFmt
$ws.Cells.Item(67,2).Value = "0x01,0x90"
$ws.Cells.Item(67,3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
$ws.Cells.Item(67,4).Value = "0x01,0x90,"
$ws.Cells.Item(67,5).Value = "0x3"
$ws.Cells.Item(67,6).Value = 400
$ws.Cells.Item(67,7).Value = $g
$ws.Cells.Item(67,8).Value = 400
$ws.Cells.Item(67,9).Value = 3
$ws.Cells.Item(68,1).Value = 45722.22834673611
$ws.Cells.Item(68,1).NumberFormat = $dateFmt
$ws.Cells.Item(68,2).Value = "0x01,0x90"
$ws.Cells.Item(68,3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
$ws.Cells.Item(68,4).Value = "0x01,0x90,"
$ws.Cells.Item(68,5).Value = "0x3"
$ws.Cells.Item(68,6).Value = 400
$ws.Cells.Item(68,7).Value = $g
$ws.Cells.Item(68,8).Value = 400
$ws.Cells.Item(68,9).Value = 3
$ws.Cells.Item(69,1).Value = 45722.22836789352
$ws.Cells.Item(69,1).NumberFormat = $dateFmt
$ws.Cells.Item(69,2).Value = "0x01,0x90"
$ws.Cells.Item(69,3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
$ws.Cells.Item(69,4).Value = "0x01,0x90,"
$ws.Cells.Item(69,5).Value = "0x3"
$ws.Cells.Item(69,6).Value = 400
$ws.Cells.Item(69,7).Value = $g
$ws.Cells.Item(69,8).Value = 400
$ws.Cells.Item(69,9).Value = 3
$ws.Cells.Item(70,1).Value = 45722.22839168982
$ws.Cells.Item(70,1).NumberFormat = $dateFmt
$ws.Cells.Item(70,2).Value = "0x01,0x90"
$ws.Cells.Item(70,3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
$ws.Cells.Item(70,4).Value = "0x01,0x90,"
$ws.Cells.Item(70,5).Value = "0x3"
$ws.Cells.Item(70,6).Value = 400
$ws.Cells.Item(70,7).Value = $g
$ws.Cells.Item(70,8).Value = 400
$ws.Cells.Item(70,9).Value = 3
$ws.Cells.Item(71,1).Value = 45722.72848770834
$ws.Cells.Item(71,1).NumberFormat = $dateFmt
$ws.Cells.Item(71,2).Value = "0x01,0x90"
$ws.Cells.Item(71,3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
$ws.Cells.Item(71,4).Value = "0x01,0x90,"
$ws.Cells.Item(71,5).Value = "0x3"
$ws.Cells.Item(71,6).Value = 400
$ws.Cells.Item(71,7).Value = $g
$ws.Cells.Item(71,8).Value = 400
$ws.Cells.Item(71,9).Value = 3
$ws.Cells.Item(72,1).Value = 45722.72850997685
$ws.Cells.Item(72,1).NumberFormat = $dateFmt
$ws.Cells.Item(72,2).Value = "0x01,0x90"
$ws.Cells.Item(72,3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
$ws.Cells.Item(72,4).Value = "0x01,0x90,"
$ws.Cells.Item(72,5).Value = "0x3"
$ws.Cells.Item(72,6).Value = 400
$ws.Cells.Item(72,7).Value = $g
$ws.Cells.Item(72,8).Value = 400
$ws.Cells.Item(72,9).Value = 3
$ws.Cells.Item(73,1).Value = 45722.72853335648
$ws.Cells.Item(73,1).NumberFormat = $dateFmt
$ws.Cells.Item(73,2).Value = "0x01,0x90"
$ws.Cells.Item(73,3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
$ws.Cells.Item(73,4).Value = "0x01,0x90,"
$ws.Cells.Item(73,5).Value = "0x3"
$ws.Cells.Item(73,6).Value = 400
$ws.Cells.Item(73,7).Value = $g
$ws.Cells.Item(73,8).Value = 400
$ws.Cells.Item(73,9).Value = 3
$ws.Cells.Item(74,1).Value = 45723.22863142361
$ws.Cells.Item(74,1).NumberFormat = $dateFmt
$ws.Cells.Item(74,2).Value = "0x01,0x90"
$ws.Cells.Item(74,3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
$ws.Cells.Item(74,4).Value = "0x01,0x90,"
$ws.Cells.Item(74,5).Value = "0x3"
$ws.Cells.Item(74,6).Value = 400
$ws.Cells.Item(74,7).Value = $g
$ws.Cells.Item(74,8).Value = 400
$ws.Cells.Item(74,9).Value = 3
$ws.Cells.Item(75,1).Value = 45723.22865329861
$ws.Cells.Item(75,1).NumberFormat = $dateFmt
$ws.Cells.Item(75,2).Value = "0x01,0x90"
$ws.Cells.Item(75,3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
$ws.Cells.Item(75,4).Value = "0x01,0x90,"
$ws.Cells.Item(75,5).Value = "0x3"
$ws.Cells.Item(75,6).Value = 400
$ws.Cells.Item(75,7).Value = $g
$ws.Cells.Item(75,8).Value = 400
$ws.Cells.Item(75,9).Value = 3
$ws.Cells.Item(76,1).Value = 45723.22867643519
$ws.Cells.Item(76,1).NumberFormat = $dateFmt
$ws.Cells.Item(76,2).Value = "0x01,0x90"
$ws.Cells.Item(76,3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
$ws.Cells.Item(76,4).Value = "0x01,0x90,"
$ws.Cells.Item(76,5).Value = "0x3"
$ws.Cells.Item(76,6).Value = 400
$ws.Cells.Item(76,7).Value = $g
$ws.Cells.Item(76,8).Value = 400
$ws.Cells.Item(76,9).Value = 3
